# Adding different messages to the InputException for different input types
# Append three new screening-log rows (10, 11, 12) to the Screening_Log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        $Worksheet,
        [string]$Cell,
        [string]$Value
    )
    # Force the cell to remain a text value (avoids Excel auto-converting
    # strings that look like dates/numbers, e.g. "2017-10-28" or "9"),
    # then restore the default "General" number format so no residual
    # explicit style is left behind on the cell.
    $rng = $Worksheet.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.NumberFormat = "General"
}

# ---- Row 10 ----
Set-TextValue $ws "A10" "2017-10-28"
Set-TextValue $ws "B10" "14:08:37.368626"
Set-TextValue $ws "C10" "k"
Set-TextValue $ws "D10" "k"
Set-TextValue $ws "E10" "9"
Set-TextValue $ws "F10" "m"
Set-TextValue $ws "G10" "y"
Set-TextValue $ws "H10" "n"
Set-TextValue $ws "I10" "n"
Set-TextValue $ws "J10" "u"
Set-TextValue $ws "K10" "j"

# ---- Row 11 ----
Set-TextValue $ws "A11" "2017-10-28"
Set-TextValue $ws "B11" "22:01:06.972138"
Set-TextValue $ws "C11" "k"
Set-TextValue $ws "D11" "k"
Set-TextValue $ws "E11" "9"
Set-TextValue $ws "F11" "m"
Set-TextValue $ws "G11" "n"
Set-TextValue $ws "H11" "k"
Set-TextValue $ws "I11" "n"
Set-TextValue $ws "J11" "k"
Set-TextValue $ws "K11" "k"

# ---- Row 12 ----
Set-TextValue $ws "A12" "2017-10-28"
Set-TextValue $ws "B12" "22:03:31.099802"
Set-TextValue $ws "C12" "k"
Set-TextValue $ws "D12" "k"
Set-TextValue $ws "E12" "9"
Set-TextValue $ws "F12" "m"
Set-TextValue $ws "G12" "n"
Set-TextValue $ws "H12" "j"
Set-TextValue $ws "I12" "n"
Set-TextValue $ws "J12" "n"
Set-TextValue $ws "K12" "n"

Write-Host "Added rows 10-12 to $($ws.Name)"
